# Apply cryptos list price/volume updates per commit "Updated cryptos list on Sat Mar 18 14:23:32 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.597.43"
$ws.Range("E2").Value = "  +3.86%  "
$ws.Range("D3").Value = "1.824.92"
$ws.Range("E3").Value = "  +4.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "341.56"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3827"
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3536"
$ws.Range("E8").Value = "  +4.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.81"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.240"
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07746"
$ws.Range("E11").Value = "  +3.93%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.41"
$ws.Range("E13").Value = "  +9.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.628"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").Value = "1.824.50"
$ws.Range("E15").Value = "  +4.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.208"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001129"
$ws.Range("E17").Value = "  +4.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06725"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.30"
$ws.Range("E19").Value = "  +4.59%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.70"
$ws.Range("E21").Value = "  +5.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.548"
$ws.Range("E22").Value = "  +5.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.17"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").Value = "27.590.73"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.480"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.680"
$ws.Range("E26").Value = "  +9.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.15"
$ws.Range("E27").Value = "  +12.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.483"
$ws.Range("E28").Value = "  +5.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.15"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "2.031.10"
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.65"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.350"
$ws.Range("E32").Value = "  +4.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.085"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.97"
$ws.Range("E34").Value = "  +7.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08810"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.656"
$ws.Range("E37").Value = "  +4.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7043"
$ws.Range("E38").Value = "  +13.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.162"
$ws.Range("E39").Value = "  +6.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2267"
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02411"
$ws.Range("E41").Value = "  +2.55%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.06504"
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.296"
$ws.Range("E43").Value = "  +5.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.75"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6632"
$ws.Range("E45").Value = "  +9.77%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.935"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.189"
$ws.Range("E48").Value = "  +6.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "133.51"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07322"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.46"
$ws.Range("E51").Value = "  +4.77%  "
